# Generate Report for Handoff
# Replace the old GUID-based file identifiers with the new ones, update the
# handoff/handback timestamps, and refresh hyperlink display text to match.

$wb = $excel.ActiveWorkbook

$oldGuid = "7f49b3e0-872f-49c7-bfef-baa7b1f949b3"
$newGuid = "84b9a225-6938-4a42-9132-98900a48dad2"

$oldHash = "5b0d3032478dde07b833b1999e425d8643b2af29"
$newHash = "29f528086287d2f22261ee3ad60e35b5ebb050a2"

$hyperlinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/50506367406b363e508fc014c74d43d2f685e00a/e2e/$oldGuid.md"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid.md"

$newDisplayB2 = "e2e\$newGuid.md"
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $hyperlinkAddress, [Type]::Missing, [Type]::Missing, $newDisplayB2)

$wsOverview.Range("G2").Value = "2016-08-18 19:01:35"

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$newDisplayA2 = "$newGuid.md"
$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $hyperlinkAddress, [Type]::Missing, [Type]::Missing, $newDisplayA2)

$wsZhCn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-18 19:01:30"

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $hyperlinkAddress, [Type]::Missing, [Type]::Missing, $newDisplayA2)

$wsDeDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-18 19:01:35"
